$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the board-state strings referenced by A1:A8 (and the derived C/E formula
# cells recalculate automatically), plus the trailing CLI-argument string in E6.
# The order below matches the order the shared-string table ends up in.
$ws.Range("A2").Value = " 1, 1, 1, 0, 0, 1, 1, 1,"
$ws.Range("A8").Value = "-5,-2,-3,-8,-9,-3,-2,-5"
$ws.Range("E6").Value = "`" 0 1"
$ws.Range("A4").Value = " 0,-1, 0, 1, 0, 3, 0, 0,"
$ws.Range("A1").Value = "5, 2, 0, 8, 9, 0, 0, 5,"
$ws.Range("A3").Value = " 0, 0, 0, 3, 1, 2, 0, 0,"
$ws.Range("A6").Value = " 0, 0,-1, 0, 0, 0,-1, 0,"
$ws.Range("A7").Value = " 0, 0, 0,-1,-1, 0, 0,-1,"
$ws.Range("A5").Value = " -1, 0, 0, 0, 0, 0, 0, 0,"
